$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "Лист1"
$ws2.Name = "Расчет"

# --- Sheet1: zoom change (closest achievable via COM: set window zoom while sheet is active) ---
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100

# --- Sheet2: view / selection + column width ---
$ws2.Activate()
$ws2.Columns.Item(1).ColumnWidth = 27

# Convert O2:O27 / Q2:Q27 into shared formulas (matches original formula text)
$ws2.Range("O2:O27").Formula = "=L2-N2"
$ws2.Range("Q2:Q27").Formula = "=ABS(O2)/L2"

# Update selection to match the authored edit
$ws2.Range("L100").Select()
